# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
# - Shrink the "Status" column(s) to match the new, narrower report width

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColWidth = 13.4101845877511

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($addr in @("C2", "C3")) {
    $cell = $wsZhCn.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($addr in @("C2", "C3")) {
    $cell = $wsDeDe.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
